$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("z2,1")

# Row 1: shift the "1" marker from C1 to F1
$ws.Range("C1").Value = 0
$ws.Range("D1").Value = 0
$ws.Range("E1").Value = 0
$ws.Range("F1").Value = 1

# Row 3: clear the marker at J3
$ws.Range("J3").Value = 0

# Row 5: set marker at H5
$ws.Range("H5").Value = 1

# Row 6: set marker at E6
$ws.Range("E6").Value = 1

# Row 7: set marker at J7
$ws.Range("J7").Value = 1

# Row 8: set marker at G8
$ws.Range("G8").Value = 1
